# Ticket 002: nice header and background created for both pages
# - Add PR link (F2) and mirror the "done" status into G3
# - Tidy up alignment/formatting used across the sheet (consolidates a
#   handful of duplicate cell styles down to the ones Excel kept)
# - Widen column F so the PR link is readable, and let row heights settle
#   back down now that several descriptions wrap onto fewer lines

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content -----------------------------------------------------
$ws.Range("F2").Value = "https://github.com/zsigmondO/rick-and-morty-catalog/pull/1"
$ws.Range("G3").Value = "done"

# --- Header row (A1:B1): keep the red/bold header look, just re-assert
#     the centered alignment so it lines up with the rest of the header
$headerAB = $ws.Range("A1:B1")
$headerAB.HorizontalAlignment = -4108
$headerAB.VerticalAlignment = -4107

# --- Ticket-number / branch-name columns (A/B) plus the two "Status"
#     cells in column G: centered horizontally and vertically, normal
#     (non-bold) font
$centered = $ws.Range("A2,B2,G2,A3,B3,G3,A5,B5,A6,B6,A7,B7,A8,B8,A9,B9,A11,B11,A12,B12,A13,B13,A15,B15,A16,B16,A17,B17,A18,B18,A19,B19")
$centered.HorizontalAlignment = -4108
$centered.VerticalAlignment = -4108
$centered.Font.Bold = $false

# --- Section header row (A4:F4 "Home page tickets"): bold, centered
#     horizontally only
$section1 = $ws.Range("A4:F4")
$section1.HorizontalAlignment = -4108
$section1.VerticalAlignment = -4107
$section1.Font.Bold = $true

# --- C7 description cell: left aligned, vertical centered, wrapped
$c7 = $ws.Range("C7")
$c7.HorizontalAlignment = -4131
$c7.VerticalAlignment = -4108
$c7.WrapText = $true

# --- Section header rows (A10:E10 "Profile page tickets" and A14:E14
#     "Misc. or optional tickets"): bold, centered horizontally and
#     vertically
$section2 = $ws.Range("A10:E10,A14:E14")
$section2.HorizontalAlignment = -4108
$section2.VerticalAlignment = -4108
$section2.Font.Bold = $true

# --- Column widths / row heights --------------------------------------
$ws.Columns("F").ColumnWidth = 51.25

$ws.Rows(5).RowHeight = 43.2
$ws.Rows(6).RowHeight = 57.6
$ws.Rows(7).RowHeight = 43.2
$ws.Rows(8).RowHeight = 57.6
$ws.Rows(9).RowHeight = 28.8
$ws.Rows(11).RowHeight = 28.8

# --- Leave the selection where the author finished typing -------------
$ws.Range("F2").Select()
